$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 1082-1096: column Q ("backup") resolves from an empty
#     placeholder (inlineStr) to a numeric value. ---
$qUpdates = @{
    1082 = 2
    1083 = 0
    1084 = 0
    1085 = 0
    1086 = 0
    1087 = 0
    1088 = 0
    1089 = 0
    1090 = 0
    1091 = 0
    1092 = 0
    1093 = 0
    1094 = 0
    1095 = 0
    1096 = 0
}
foreach ($r in $qUpdates.Keys) {
    $ws.Cells.Item($r, 17).Value = $qUpdates[$r]
}

# --- New rows 1097-1110: Datetime .. detect_structure (cols A-P);
#     column Q ("backup") is left as an empty placeholder, same as
#     the newest rows were before this edit. ---
$newRows = @(
    @(1097, 45680, 1065, 1088, 1059, 1082.849975585938, 1467332, 2025, 1, 23, 0, 0, 0, 4, 0, 0, 0),
    @(1098, 45681, 1080, 1087, 1051.949951171875, 1055.75, 756024, 2025, 1, 24, 0, 0, 0, 4, 0, 0, 0),
    @(1099, 45684, 1040.349975585938, 1042.949951171875, 1015.900024414062, 1029.949951171875, 702909, 2025, 1, 27, 0, 0, 0, 5, 0, 0, 0),
    @(1100, 45685, 1015, 1043.800048828125, 982, 1021.650024414062, 2125229, 2025, 1, 28, 0, 0, 0, 5, 2, 0, 0),
    @(1101, 45686, 1011, 1032.400024414062, 1000, 1025, 1074681, 2025, 1, 29, 0, 0, 0, 5, 0, 0, 0),
    @(1102, 45687, 1040.099975585938, 1058.400024414062, 998, 1043.849975585938, 3974586, 2025, 1, 30, 0, 0, 0, 5, 0, 0, 0),
    @(1103, 45688, 1046, 1069, 1029.75, 1061.25, 2955047, 2025, 1, 31, 0, 0, 0, 5, 0, 0, 0),
    @(1104, 45689, 1073.449951171875, 1117.75, 1067.25, 1108.25, 1484218, 2025, 2, 1, 0, 0, 0, 5, 0, 0, 0),
    @(1105, 45691, 1125, 1150, 1090.550048828125, 1145, 2635417, 2025, 2, 3, 0, 0, 0, 6, 0, 0, 0),
    @(1106, 45692, 1150.199951171875, 1184, 1145, 1174.300048828125, 2429222, 2025, 2, 4, 0, 0, 0, 6, 0, 0, 0),
    @(1107, 45693, 1180, 1200.949951171875, 1164.25, 1169.699951171875, 1565195, 2025, 2, 5, 0, 0, 0, 6, 0, 0, 0),
    @(1108, 45694, 1185, 1187.949951171875, 1121, 1125.75, 1003312, 2025, 2, 6, 0, 0, 0, 6, 0, 0, 0),
    @(1109, 45695, 1133, 1147.050048828125, 1110.599975585938, 1143.199951171875, 1906494, 2025, 2, 7, 0, 0, 0, 6, 0, 0, 0),
    @(1110, 45698, 1140, 1140, 1065, 1071.849975585938, 3230528, 2025, 2, 10, 0, 0, 0, 7, 0, 0, 0),
)

foreach ($row in $newRows) {
    $rowNum = $row[0]
    for ($col = 1; $col -le 16; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $row[$col]
    }
    # Column A carries the same date/time display format as the rest of the sheet.
    $ws.Cells.Item($rowNum, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
